$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8304660320281982
$ws.Range("B1").Value = 3.695141315460205
$ws.Range("C1").Value = 6.106452465057373
$ws.Range("D1").Value = 2.738850355148315
$ws.Range("E1").Value = 1.903808832168579
